$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '22.357.20'
$ws.Range('E2').Value = '  -0.02%  '
$ws.Range('D3').Value = '1.565.85'
$ws.Range('E3').Value = '  +0.12%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.005'
$ws.Range('E4').Value = '  +0.29%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '1.005'
$ws.Range('E5').Value = '  +0.35%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '289.80'
$ws.Range('E6').Value = '  +0.03%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.3749'
$ws.Range('E7').Value = '  +0.99%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '49.23'
$ws.Range('E8').Value = '  +0.28%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.3373'
$ws.Range('E9').Value = '  -0.78%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.07507'
$ws.Range('E10').Value = '  -1.81%  '
$ws.Range('E11').Value = '  -3.87%  '
$ws.Range('E12').Value = '  +0.35%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '20.81'
$ws.Range('E13').Value = '  -2.76%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '5.895'
$ws.Range('E14').Value = '  -2.57%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '6.865'
$ws.Range('E15').Value = '  -0.84%  '
$ws.Range('D16').Value = '1.566.86'
$ws.Range('E16').Value = '  +0.58%  '
$ws.Range('E17').Value = '  -1.18%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '89.26'
$ws.Range('E18').Value = '  -0.88%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.06714'
$ws.Range('E19').Value = '  -0.20%  '
$ws.Range('E20').Value = '  +0.39%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '6.169'
$ws.Range('E21').Value = '  -1.15%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '16.31'
$ws.Range('E22').Value = '  -1.46%  '
$ws.Range('E23').Value = '  -1.74%  '
$ws.Range('D24').Value = '22.367.31'
$ws.Range('E24').Value = '  +0.07%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.369'
$ws.Range('E25').Value = '  -0.84%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '2.687'
$ws.Range('E26').Value = '  -4.73%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '19.97'
$ws.Range('E27').Value = '  -0.96%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '147.67'
$ws.Range('E28').Value = '  +1.72%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '4.996'
$ws.Range('E29').Value = '  +0.31%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '124.86'
$ws.Range('E30').Value = '  -0.31%  '
$ws.Range('D31').Value = '1.740.03'
$ws.Range('E31').Value = '  +0.36%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '2.017'
$ws.Range('E32').Value = '  +0.28%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.9843'
$ws.Range('E33').Value = '  -2.20%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '5.941'
$ws.Range('E34').Value = '  -4.28%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '9.814'
$ws.Range('E35').Value = '  -2.18%  '
$ws.Range('B36').Value = 'Stellar'
$ws.Range('C36').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.08448'
$ws.Range('E36').Value = '  +0.08%  '
$ws.Range('B37').Value = 'TrustWalletToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '1.404'
$ws.Range('E37').Value = '  +9.47%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.02455'
$ws.Range('E38').Value = '  -3.11%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.2260'
$ws.Range('E39').Value = '  -2.72%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.06427'
$ws.Range('E40').Value = '  +0.27%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '5.359'
$ws.Range('E41').Value = '  -3.06%  '
$ws.Range('E42').Value = '  -1.61%  '
$ws.Range('E43').Value = '  -6.32%  '
$ws.Range('E44').Value = '  +0.38%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '13.79'
$ws.Range('E45').Value = '  -2.19%  '
$ws.Range('E46').Value = '  +0.76%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.5792'
$ws.Range('E47').Value = '  -3.04%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '2.046'
$ws.Range('E48').Value = '  -2.30%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.249'
$ws.Range('E49').Value = '  -1.46%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '123.89'
$ws.Range('E50').Value = '  -0.70%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.07300'
$ws.Range('E51').Value = '  +0.47%  '
